# Reformat the single question-bank JSON blob in A-column from a compact
# Python-literal one-liner into pretty-printed JSON (matching the authors
# local re-export), drop the old header/flag row, and drop the now-unused
# bold+bordered cell style so the remaining cell is plain text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 held the old flag cell (styled, numeric 0); remove it so the question
# text (old A2) shifts up to become A1, unstyled.
$ws.Rows.Item(1).Delete()

# Stage the new text on a scratch cell first, then Cut/Paste it into A1.
# Writing multi-line text straight into A1 would leave Excel's automatic
# row-autofit height stamped on row 1; cutting the already-sized scratch
# cell into place moves the text without carrying that row metadata.
$scratch = $ws.Range("ZZ999")
$scratch.Value = 'questions = [
    {
        "title": "Review the following code structure:class WorkingCenter(models.Model):\n  name = CharField(max_length=100)\n\nclass Employee(models.Model):\n  name = CharField(max_length=100)\n  last_name = CharField(max_length=100)  \n  working_center = models.ForeignKey(WorkingCenter)\n\nclass Contract(models.Model):\n  start_date = models.DateTimeField()\n  end_date = models.DateTimeField(null=True)\n  employee = models.ForeignKey(Employee)\nCan I execute the following code to access the list of contracts of an employee?employee = Employee.objects.get(id=1)\ncontracts = employee.contracts.all()",
        "ques_type": 2,
        "options": [
            "Yes, and the result will be a queryset.",
            "It is not possible to access the contract list directly from an employee. To do so, it is necessary to execute the following line of code:contracts = Contract.objects.filter(employee=employee) \n",
            "The code will fail, but changing the Contract model will make it work, as follows:employee = models.ForeignKey(Employee,  on_delete=models.CASCADE, related_name=''contracts'')\n",
            "No, the code will fail. The correct line of code is:contracts = employee.get_contracts.all() \nThe reason is that every Model implements generic getters to access linked models."
        ],
        "score": "The code will fail, but changing the Contract model will make it work, as follows:employee = models.ForeignKey(Employee,  on_delete=models.CASCADE, related_name=''contracts'')"
    },
    {
        "title": "Select the existing fields for a ModelSerializer definition.",
        "ques_type": 15,
        "options": [
            "SerializerMethodField",
            "CharField",
            "BooleanField",
            "IntegerField",
            "StringField",
            "MatrixField",
            "DocumentField",
            "ObjectField"
        ],
        "score": [
            "SerializerMethodField",
            "CharField",
            "BooleanField",
            "IntegerField"
        ]
    },
    {
        "title": "Which methods could you override to change the behavior of a ModelViewSet view that creates a new object?",
        "ques_type": 15,
        "options": [
            "def perform_create(self, serializer): \n",
            "def create(self, request, *args, **kwargs): \n",
            "def post(self, request, *args, **kwargs): \n",
            "def put(self, request, *args, **kwargs):\n"
        ],
        "score": [
            "def perform_create(self, serializer):",
            "def create(self, request, *args, **kwargs):"
        ]
    },
    {
        "title": "Using the Django REST Framework, you are developing an API in which you have one endpoint to create an employee. To do so, your endpoint will expect the following JSON structure:{\n  name: \"Employee name\",\n  last_name: \"Employee last name\",\n  expertise: \"Description of the expertise\"\n}\nWhich component of DRF is responsible for validating every input field?",
        "ques_type": 2,
        "options": [
            "The serializer\t",
            "The view",
            "The model",
            "The validation middleware"
        ],
        "score": "The serializer"
    }
]'
$scratch.Cut($ws.Range("A1"))

# Clean up the now-empty scratch row so no stray formatting lingers.
$ws.Rows.Item(999).Delete()

